$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data.
# Cells whose new text would be auto-parsed as a number (e.g. "0.9976", "1.000")
# are forced to remain plain text by temporarily formatting them as Text,
# then resetting the cell style back to Normal so no numeric conversion/
# extra styling is left behind - matching the original plain inline-string cells.

$ws.Range("D2").Value = "26.137.50"
$ws.Range("E2").Value = "  +4.43%  "
$ws.Range("D3").Value = "1.690.94"
$ws.Range("E3").Value = "  +3.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9976"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9986"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4674"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2627"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06183"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.96%  "
$ws.Range("D10").Value = "1.683.20"
$ws.Range("E10").Value = "  +2.60%  "
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.407"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5848"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "75.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9986"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9986"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "26.121.81"
$ws.Range("E18").Value = "  +4.47%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006767"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.15%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.87%  "
$ws.Range("D21").Value = "1.897.23"
$ws.Range("E21").Value = "  +2.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.538"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.758"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.286"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "134.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.74%  "
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.739"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.987"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.675"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07759"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04376"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.596"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("E35").Value = "  +4.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9634"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9250"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "110.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.387"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.03%  "
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.889"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01467"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3764"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.077"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1134"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.23%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05319"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.44%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.211"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.684"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.46%  "
$ws.Range("E50").Value = "  +2.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.000"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.04%  "
